$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 49869.168
$ws.Range("I10").Value = 500
$ws.Range("J10").Value = 59743
$ws.Range("K10").Value = 500
$ws.Range("L10").Value = 59743
$ws.Range("M10").Value = -207
$ws.Range("N10").Value = -60329

$ws.Range("H107").Value = 893.54285
$ws.Range("I107").Value = 930.1724
$ws.Range("J107").Value = 716.5
$ws.Range("K107").Value = 930.1724
$ws.Range("L107").Value = 716.5
$ws.Range("M107").Value = 989.8276
$ws.Range("N107").Value = -4556.5

$ws.Range("H116").Value = 17957418
$ws.Range("I116").Value = 13335694
$ws.Range("J116").Value = 20845998
$ws.Range("K116").Value = 13335694
$ws.Range("L116").Value = 20845998
$ws.Range("M116").Value = -13332252
$ws.Range("N116").Value = -20852882

$ws.Range("H132").Value = 4447685
$ws.Range("I132").Value = 3435.4348
$ws.Range("J132").Value = 55556556
$ws.Range("K132").Value = 10306.3044
$ws.Range("L132").Value = 166669668
$ws.Range("M132").Value = -7776.304400000001
$ws.Range("N132").Value = -166674728

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 2501002.5
$ws.Range("I8").Value = 2501002.5
$ws.Range("K8").Value = 2501002.5
$ws.Range("M8").Value = -2500858.5

$ws.Range("H45").Value = 1411.625
$ws.Range("I45").Value = 1209.1818
$ws.Range("J45").Value = 1857
$ws.Range("K45").Value = 1209.1818
$ws.Range("L45").Value = 1857
$ws.Range("M45").Value = -832.1818000000001
$ws.Range("N45").Value = -2611

$ws.Range("H63").Value = 1486.3334
$ws.Range("I63").Value = 1486.6666
$ws.Range("J63").Value = 1485.3334
$ws.Range("K63").Value = 1486.6666
$ws.Range("L63").Value = 1485.3334
$ws.Range("M63").Value = -800.6666
$ws.Range("N63").Value = -2857.3334

$ws.Range("H66").Value = 1486.3334
$ws.Range("I66").Value = 1486.6666
$ws.Range("J66").Value = 1485.3334
$ws.Range("K66").Value = 7433.333000000001
$ws.Range("L66").Value = 7426.666999999999
$ws.Range("M66").Value = -4001.333000000001
$ws.Range("N66").Value = -14290.667

$ws.Range("H74").Value = 52875144
$ws.Range("I74").Value = 45455856
$ws.Range("K74").Value = 45455856
$ws.Range("M74").Value = -45454982

$ws.Range("H77").Value = 52875144
$ws.Range("I77").Value = 45455856
$ws.Range("K77").Value = 227279280
$ws.Range("M77").Value = -227274912

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 497.57144
$ws.Range("I12").Value = 413.83334
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 413.83334
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = -245.83334
$ws.Range("N12").Value = -1336

$ws.Range("H14").Value = 5000
$ws.Range("J14").Value = 5000
$ws.Range("L14").Value = 5000
$ws.Range("N14").Value = -5344

$ws.Range("H20").Value = 26330656
$ws.Range("I20").Value = 38470852
$ws.Range("J20").Value = 26891.334
$ws.Range("K20").Value = 38470852
$ws.Range("L20").Value = 26891.334
$ws.Range("M20").Value = -38470605
$ws.Range("N20").Value = -27385.334

$ws.Range("H80").Value = 222.14285
$ws.Range("I80").Value = 171
$ws.Range("K80").Value = 171
$ws.Range("M80").Value = 827

$ws.Range("H83").Value = 222.14285
$ws.Range("I83").Value = 171
$ws.Range("K83").Value = 855
$ws.Range("M83").Value = 4137

$ws.Range("H99").Value = 1586.5
$ws.Range("I99").Value = 1297.125
$ws.Range("J99").Value = 2049.5
$ws.Range("K99").Value = 1297.125
$ws.Range("L99").Value = 2049.5
$ws.Range("M99").Value = 200.875
$ws.Range("N99").Value = -5045.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 271.75
$ws.Range("I5").Value = 82.625
$ws.Range("K5").Value = 82.625
$ws.Range("M5").Value = 29.375

$ws.Range("H11").Value = 300
$ws.Range("I11").Value = 300
$ws.Range("K11").Value = 300
$ws.Range("M11").Value = -160

$ws.Range("H12").Value = 30741.3
$ws.Range("I12").Value = 279
$ws.Range("J12").Value = 61203.6
$ws.Range("K12").Value = 279
$ws.Range("L12").Value = 61203.6
$ws.Range("M12").Value = -109
$ws.Range("N12").Value = -61543.6

$ws.Range("H31").Value = 2318542.5
$ws.Range("I31").Value = 2779316.5
$ws.Range("J31").Value = 1899657
$ws.Range("K31").Value = 2779316.5
$ws.Range("L31").Value = 1899657
$ws.Range("M31").Value = -2779021.5
$ws.Range("N31").Value = -1900247

$ws.Range("H34").Value = 2318542.5
$ws.Range("I34").Value = 2779316.5
$ws.Range("J34").Value = 1899657
$ws.Range("K34").Value = 2779316.5
$ws.Range("L34").Value = 1899657
$ws.Range("M34").Value = -2779114.5
$ws.Range("N34").Value = -1900061

$ws.Range("H122").Value = 8547.817999999999
$ws.Range("I122").Value = 8547.817999999999
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 25643.454
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -23193.454
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 2068.2
$ws.Range("I132").Value = 1325.9474
$ws.Range("K132").Value = 3977.8422
$ws.Range("M132").Value = -1447.8422

$ws.Range("H134").Value = 1606922.4
$ws.Range("I134").Value = 6787.25
$ws.Range("K134").Value = 20361.75
$ws.Range("M134").Value = -17826.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1488.421
$ws.Range("I114").Value = 202.66667
$ws.Range("J114").Value = 1729.5
$ws.Range("K114").Value = 608.00001
$ws.Range("L114").Value = 5188.5
$ws.Range("M114").Value = 2645.99999
$ws.Range("N114").Value = -11696.5

$ws.Range("H132").Value = 1841.3889
$ws.Range("I132").Value = 914.3
$ws.Range("J132").Value = 3000.25
$ws.Range("K132").Value = 8228.699999999999
$ws.Range("L132").Value = 27002.25
$ws.Range("M132").Value = -5698.699999999999
$ws.Range("N132").Value = -32062.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H46").Value = 9264.223
$ws.Range("J46").Value = 9922.25
$ws.Range("L46").Value = 9922.25
$ws.Range("N46").Value = -10234.25

$ws.Range("H80").Value = 17460.77
$ws.Range("I80").Value = 5998.75
$ws.Range("J80").Value = 35800
$ws.Range("K80").Value = 5998.75
$ws.Range("L80").Value = 35800
$ws.Range("M80").Value = -5000.75
$ws.Range("N80").Value = -37796

$ws.Range("H83").Value = 17460.77
$ws.Range("I83").Value = 5998.75
$ws.Range("J83").Value = 35800
$ws.Range("K83").Value = 29993.75
$ws.Range("L83").Value = 179000
$ws.Range("M83").Value = -25001.75
$ws.Range("N83").Value = -188984

$ws.Range("H126").Value = 4719.636
$ws.Range("I126").Value = 6150
$ws.Range("J126").Value = 2216.5
$ws.Range("K126").Value = 18450
$ws.Range("L126").Value = 6649.5
$ws.Range("M126").Value = -15980
$ws.Range("N126").Value = -11589.5

$ws.Range("H132").Value = 26044472
$ws.Range("I132").Value = 41274380
$ws.Range("J132").Value = 12990266
$ws.Range("K132").Value = 123823140
$ws.Range("L132").Value = 38970798
$ws.Range("M132").Value = -123820610
$ws.Range("N132").Value = -38975858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 500
$ws.Range("I3").Value = 500
$ws.Range("K3").Value = 500
$ws.Range("M3").Value = -388

$ws.Range("H15").Value = 500
$ws.Range("I15").Value = 500
$ws.Range("K15").Value = 500
$ws.Range("M15").Value = -330

$ws.Range("H46").Value = 694.80646
$ws.Range("I46").Value = 617.9666999999999
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 617.9666999999999
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -429.9666999999999
$ws.Range("N46").Value = -3376

$ws.Range("H132").Value = 15877048
$ws.Range("I132").Value = 35715484
$ws.Range("K132").Value = 107146452
$ws.Range("M132").Value = -107143922

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 10666.667
$ws.Range("J8").Value = 10000
$ws.Range("L8").Value = 10000
$ws.Range("N8").Value = -10280

$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H12").Value = 70007
$ws.Range("J12").Value = 70007
$ws.Range("L12").Value = 70007
$ws.Range("N12").Value = -70291

$ws.Range("H13").Value = 105
$ws.Range("I13").Value = 105
$ws.Range("K13").Value = 105
$ws.Range("M13").Value = 35

$ws.Range("H126").Value = 5682592.5
$ws.Range("I126").Value = 8621352
$ws.Range("J126").Value = 990.4666999999999
$ws.Range("K126").Value = 25864056
$ws.Range("L126").Value = 2971.4001
$ws.Range("M126").Value = -25861586
$ws.Range("N126").Value = -7911.4001

$ws.Range("H136").Value = 3868.8845
$ws.Range("I136").Value = 1238.1333
$ws.Range("J136").Value = 4935.4053
$ws.Range("K136").Value = 3714.3999
$ws.Range("L136").Value = 14806.2159
$ws.Range("M136").Value = -1164.3999
